$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B63").Value = 6832705
$ws.Range("B64").Value = 6832544
$ws.Range("F63").Value = "Al Bataeh"
$ws.Range("F64").Value = "Al Jazira SC"
$ws.Range("G63").Value = "Baniyas SC"
$ws.Range("G64").Value = "Ajman SCC"
$ws.Range("H63").Value = 2
$ws.Range("H64").Value = 1
$ws.Range("I63").Value = 1
$ws.Range("I64").Value = 5
$ws.Range("J63").Value = "H"
$ws.Range("J64").Value = "A"
$ws.Range("K63").Value = 2.6
$ws.Range("K64").Value = 1.222
$ws.Range("L63").Value = 3.4
$ws.Range("L64").Value = 6.5
$ws.Range("M63").Value = 2.6
$ws.Range("M64").Value = 11
$ws.Range("N63").Value = 2.25
$ws.Range("N64").Value = 1.4
$ws.Range("O63").Value = 3.4
$ws.Range("O64").Value = 5
$ws.Range("P63").Value = 3.1
$ws.Range("P64").Value = 6.5
$ws.Range("Q63").Value = -0.25
$ws.Range("Q64").Value = -1.5
$ws.Range("R63").Value = 2
$ws.Range("R64").Value = 1.95
$ws.Range("S63").Value = 1.8
$ws.Range("S64").Value = 1.85
$ws.Range("T63").Value = 2.75
$ws.Range("T64").Value = 3.75
$ws.Range("U63").Value = 1.85
$ws.Range("U64").Value = 1.95
$ws.Range("V63").Value = 1.95
$ws.Range("V64").Value = 1.85
$ws.Range("W63").Value = 1.25
$ws.Range("W64").Value = -1
$ws.Range("Y63").Value = -1
$ws.Range("Y64").Value = 5.5
$ws.Range("Z63").Value = 1
$ws.Range("Z64").Value = -1
$ws.Range("AA63").Value = -1
$ws.Range("AA64").Value = 0.8500000000000001
$ws.Range("AB63").Value = 0.425
$ws.Range("AB64").Value = 0.95
$ws.Range("AC63").Value = -0.5
$ws.Range("AC64").Value = -1
$ws.Range("B89").Value = 6832568
$ws.Range("B90").Value = 6832569
$ws.Range("F89").Value = "Hatta Dubai"
$ws.Range("F90").Value = "Emirates Club RAK"
$ws.Range("G89").Value = "Khor Fakkan"
$ws.Range("G90").Value = "Al Wasl SC"
$ws.Range("K89").Value = 2.9
$ws.Range("K90").Value = 7
$ws.Range("L89").Value = 3.6
$ws.Range("L90").Value = 6.5
$ws.Range("M89").Value = 2.15
$ws.Range("M90").Value = 1.285
$ws.Range("N89").Value = 3.1
$ws.Range("N90").Value = 10
$ws.Range("O89").Value = 3.5
$ws.Range("O90").Value = 9
$ws.Range("P89").Value = 2.05
$ws.Range("P90").Value = 1.142
$ws.Range("Q89").Value = 0.25
$ws.Range("Q90").Value = 2.5
$ws.Range("R89").Value = 1.95
$ws.Range("R90").Value = 1.8
$ws.Range("S89").Value = 1.85
$ws.Range("S90").Value = 2
$ws.Range("T89").Value = 2.75
$ws.Range("T90").Value = 4
$ws.Range("U89").Value = 1.8
$ws.Range("U90").Value = 1.95
$ws.Range("V89").Value = 2
$ws.Range("V90").Value = 1.85
$ws.Range("Y89").Value = 1.05
$ws.Range("Y90").Value = 0.1419999999999999
$ws.Range("Z89").Value = -1
$ws.Range("Z90").Value = 0.8
$ws.Range("AA89").Value = 0.8500000000000001
$ws.Range("AA90").Value = -1
$ws.Range("AC89").Value = 1
$ws.Range("AC90").Value = 0.8500000000000001
$ws.Range("B110").Value = 6832714
$ws.Range("B111").Value = 6832584
$ws.Range("F110").Value = "Emirates Club RAK"
$ws.Range("F111").Value = "Al Nasr SC"
$ws.Range("G110").Value = "Baniyas SC"
$ws.Range("G111").Value = "Hatta Dubai"
$ws.Range("I110").Value = 2
$ws.Range("I111").Value = 0
$ws.Range("J110").Value = "A"
$ws.Range("J111").Value = "H"
$ws.Range("K110").Value = 3.4
$ws.Range("K111").Value = 1.363
$ws.Range("L110").Value = 3.6
$ws.Range("L111").Value = 5
$ws.Range("M110").Value = 1.909
$ws.Range("M111").Value = 6
$ws.Range("N110").Value = 3.3
$ws.Range("N111").Value = 1.25
$ws.Range("O110").Value = 3.75
$ws.Range("O111").Value = 6
$ws.Range("P110").Value = 1.909
$ws.Range("P111").Value = 8
$ws.Range("Q110").Value = 0.5
$ws.Range("Q111").Value = -1.75
$ws.Range("U110").Value = 1.85
$ws.Range("U111").Value = 1.825
$ws.Range("V110").Value = 1.95
$ws.Range("V111").Value = 1.975
$ws.Range("W110").Value = -1
$ws.Range("W111").Value = 0.25
$ws.Range("Y110").Value = 0.909
$ws.Range("Y111").Value = -1
$ws.Range("AB110").Value = -0.5
$ws.Range("AB111").Value = -1
$ws.Range("AC110").Value = 0.475
$ws.Range("AC111").Value = 0.9750000000000001
$ws.Range("B114").Value = 6832585
$ws.Range("B116").Value = 6832715
$ws.Range("F114").Value = "Al Ain SCC"
$ws.Range("F116").Value = "Al Bataeh"
$ws.Range("G114").Value = "Sharjah SCC"
$ws.Range("G116").Value = "Al Jazira SC"
$ws.Range("H114").Value = 1
$ws.Range("H116").Value = 3
$ws.Range("I114").Value = 1
$ws.Range("I116").Value = 2
$ws.Range("J114").Value = "D"
$ws.Range("J116").Value = "H"
$ws.Range("K114").Value = 2
$ws.Range("K116").Value = 4.333
$ws.Range("L114").Value = 3.6
$ws.Range("L116").Value = 4
$ws.Range("M114").Value = 3.2
$ws.Range("M116").Value = 1.615
$ws.Range("N114").Value = 1.7
$ws.Range("N116").Value = 3.8
$ws.Range("P114").Value = 4
$ws.Range("P116").Value = 1.727
$ws.Range("Q114").Value = -0.75
$ws.Range("Q116").Value = 0.75
$ws.Range("R114").Value = 1.9
$ws.Range("R116").Value = 1.875
$ws.Range("S114").Value = 1.9
$ws.Range("S116").Value = 1.925
$ws.Range("T114").Value = 3
$ws.Range("T116").Value = 3.75
$ws.Range("U114").Value = 1.85
$ws.Range("U116").Value = 1.95
$ws.Range("V114").Value = 1.95
$ws.Range("V116").Value = 1.85
$ws.Range("W114").Value = -1
$ws.Range("W116").Value = 2.8
$ws.Range("X114").Value = 3
$ws.Range("X116").Value = -1
$ws.Range("Z114").Value = -1
$ws.Range("Z116").Value = 0.875
$ws.Range("AA114").Value = 0.8999999999999999
$ws.Range("AA116").Value = -1
$ws.Range("AB114").Value = -1
$ws.Range("AB116").Value = 0.95
$ws.Range("AC114").Value = 0.95
$ws.Range("AC116").Value = -1
$ws.Range("U117").Value = 1.925
$ws.Range("V117").Value = 1.875
$ws.Range("U122").Value = 1.975
$ws.Range("V122").Value = 1.825
